$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 619
$ws1.Range("F6").Value = 14299
$ws1.Range("F7").Value = 16362
$ws1.Range("F10").Value = 2
$ws1.Range("F19").Value = 102
$ws1.Range("F21").Value = 1249
$ws1.Range("F26").Value = 6604
$ws1.Range("F32").Value = 5715
$ws1.Range("F35").Value = 181
$ws1.Range("F36").Value = 4766

# Sheet "全部类型" (All Types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 619
$ws4.Range("F6").Value = 14299
$ws4.Range("F7").Value = 16363
$ws4.Range("F10").Value = 2
$ws4.Range("F19").Value = 102
$ws4.Range("F21").Value = 1249
$ws4.Range("F27").Value = 6605
$ws4.Range("F35").Value = 5715
$ws4.Range("F38").Value = 181
$ws4.Range("F39").Value = 4766
